$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $text) {
    $sheet.Cells.Item($row, $col).Value = "'" + $text
    $sheet.Cells.Item($row, $col).Style = "Normal"
}

function Set-NumCell($sheet, $row, $col, $num) {
    $sheet.Cells.Item($row, $col).Value = $num
}

# --- Matches_SOG: append new match rows 437-442 ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @(437, "897737", "2025-11-06T19:00:00", "Локомотив", "Спартак", 27, 29),
    @(438, "897738", "2025-11-06T19:10:00", "Динамо Мн", "Ак Барс", 34, 26),
    @(439, "897730", "2025-11-06T19:30:00", "СКА", "Драконы", 39, 25),
    @(440, "897735", "2025-11-06T19:30:00", "Динамо М", "Лада", 51, 24),
    @(441, "897736", "2025-11-06T19:30:00", "ЦСКА", "Торпедо", 29, 20),
    @(442, "897739", "2025-11-06T19:30:00", "ХК Сочи", "Северсталь", 19, 37)
)

foreach ($r in $newMatches) {
    $rowNum = $r[0]
    Set-TextCell $wsMatches $rowNum 1 $r[1]
    Set-TextCell $wsMatches $rowNum 2 $r[2]
    Set-TextCell $wsMatches $rowNum 3 $r[3]
    Set-TextCell $wsMatches $rowNum 4 $r[4]
    Set-NumCell $wsMatches $rowNum 5 $r[5]
    Set-NumCell $wsMatches $rowNum 6 $r[6]
    Set-TextCell $wsMatches $rowNum 7 "khl_text"
}

# --- Shots_HA: refresh as_of_utc + updated aggregate stats ---
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

Set-TextCell $wsShotsHA 2 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 3 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 4 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 5 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 5 6 19
Set-NumCell $wsShotsHA 5 11 630
Set-NumCell $wsShotsHA 5 12 568
Set-NumCell $wsShotsHA 5 13 33.2
Set-NumCell $wsShotsHA 5 14 29.9
Set-TextCell $wsShotsHA 6 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 7 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 8 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 8 5 16
Set-NumCell $wsShotsHA 8 7 533
Set-NumCell $wsShotsHA 8 8 439
Set-NumCell $wsShotsHA 8 9 33.3
Set-NumCell $wsShotsHA 8 10 27.4
Set-TextCell $wsShotsHA 9 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 9 5 22
Set-NumCell $wsShotsHA 9 7 787
Set-NumCell $wsShotsHA 9 8 603
Set-NumCell $wsShotsHA 9 9 35.8
Set-NumCell $wsShotsHA 9 10 27.4
Set-TextCell $wsShotsHA 10 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 10 6 20
Set-NumCell $wsShotsHA 10 11 552
Set-NumCell $wsShotsHA 10 12 737
Set-NumCell $wsShotsHA 10 13 27.6
Set-NumCell $wsShotsHA 10 14 36.9
Set-TextCell $wsShotsHA 11 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 11 6 18
Set-NumCell $wsShotsHA 11 11 481
Set-NumCell $wsShotsHA 11 12 674
Set-NumCell $wsShotsHA 11 13 26.7
Set-NumCell $wsShotsHA 11 14 37.4
Set-TextCell $wsShotsHA 12 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 12 5 17
Set-NumCell $wsShotsHA 12 7 528
Set-NumCell $wsShotsHA 12 8 466
Set-NumCell $wsShotsHA 12 9 31.1
Set-NumCell $wsShotsHA 12 10 27.4
Set-TextCell $wsShotsHA 13 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 14 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 15 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 15 5 25
Set-NumCell $wsShotsHA 15 7 830
Set-NumCell $wsShotsHA 15 8 834
Set-NumCell $wsShotsHA 15 9 33.2
Set-NumCell $wsShotsHA 15 10 33.4
Set-TextCell $wsShotsHA 16 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 17 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 17 6 23
Set-NumCell $wsShotsHA 17 11 759
Set-NumCell $wsShotsHA 17 12 599
Set-NumCell $wsShotsHA 17 13 33
Set-NumCell $wsShotsHA 17 14 26
Set-TextCell $wsShotsHA 18 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 19 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 19 6 15
Set-NumCell $wsShotsHA 19 11 518
Set-NumCell $wsShotsHA 19 12 534
Set-NumCell $wsShotsHA 19 13 34.5
Set-NumCell $wsShotsHA 19 14 35.6
Set-TextCell $wsShotsHA 20 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 20 6 27
Set-NumCell $wsShotsHA 20 11 931
Set-NumCell $wsShotsHA 20 12 866
Set-NumCell $wsShotsHA 20 13 34.5
Set-NumCell $wsShotsHA 20 14 32.1
Set-TextCell $wsShotsHA 21 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsHA 22 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 22 5 19
Set-NumCell $wsShotsHA 22 7 560
Set-NumCell $wsShotsHA 22 8 622
Set-NumCell $wsShotsHA 22 9 29.5
Set-NumCell $wsShotsHA 22 10 32.7
Set-TextCell $wsShotsHA 23 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsHA 23 5 18
Set-NumCell $wsShotsHA 23 7 412
Set-NumCell $wsShotsHA 23 8 519
Set-NumCell $wsShotsHA 23 9 22.9
Set-NumCell $wsShotsHA 23 10 28.8

# --- Shots_Summary: refresh as_of_utc + updated aggregate stats ---
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

Set-TextCell $wsShotsSummary 2 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 3 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 4 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 5 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 5 5 42
Set-NumCell $wsShotsSummary 5 6 1409
Set-NumCell $wsShotsSummary 5 7 1155
Set-NumCell $wsShotsSummary 5 8 33.5
Set-NumCell $wsShotsSummary 5 9 27.5
Set-TextCell $wsShotsSummary 6 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 7 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 8 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 8 5 37
Set-NumCell $wsShotsSummary 8 6 1129
Set-NumCell $wsShotsSummary 8 7 1112
Set-NumCell $wsShotsSummary 8 8 30.5
Set-NumCell $wsShotsSummary 8 9 30.1
Set-TextCell $wsShotsSummary 9 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 9 5 40
Set-NumCell $wsShotsSummary 9 6 1455
Set-NumCell $wsShotsSummary 9 7 1087
Set-TextCell $wsShotsSummary 10 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 10 5 39
Set-NumCell $wsShotsSummary 10 6 1084
Set-NumCell $wsShotsSummary 10 7 1399
Set-NumCell $wsShotsSummary 10 8 27.8
Set-NumCell $wsShotsSummary 10 9 35.9
Set-TextCell $wsShotsSummary 11 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 11 5 40
Set-NumCell $wsShotsSummary 11 6 1079
Set-NumCell $wsShotsSummary 11 7 1450
Set-NumCell $wsShotsSummary 11 8 27
Set-NumCell $wsShotsSummary 11 9 36.2
Set-TextCell $wsShotsSummary 12 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 12 5 42
Set-NumCell $wsShotsSummary 12 6 1305
Set-NumCell $wsShotsSummary 12 7 1074
Set-NumCell $wsShotsSummary 12 8 31.1
Set-NumCell $wsShotsSummary 12 9 25.6
Set-TextCell $wsShotsSummary 13 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 14 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 15 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 15 5 40
Set-NumCell $wsShotsSummary 15 6 1298
Set-NumCell $wsShotsSummary 15 7 1321
Set-NumCell $wsShotsSummary 15 8 32.5
Set-NumCell $wsShotsSummary 15 9 33
Set-TextCell $wsShotsSummary 16 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 17 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 17 5 39
Set-NumCell $wsShotsSummary 17 6 1228
Set-NumCell $wsShotsSummary 17 7 962
Set-NumCell $wsShotsSummary 17 8 31.5
Set-NumCell $wsShotsSummary 17 9 24.7
Set-TextCell $wsShotsSummary 18 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 19 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 19 5 39
Set-NumCell $wsShotsSummary 19 6 1372
Set-NumCell $wsShotsSummary 19 7 1198
Set-NumCell $wsShotsSummary 19 8 35.2
Set-NumCell $wsShotsSummary 19 9 30.7
Set-TextCell $wsShotsSummary 20 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 20 5 47
Set-NumCell $wsShotsSummary 20 6 1574
Set-NumCell $wsShotsSummary 20 7 1462
Set-NumCell $wsShotsSummary 20 8 33.5
Set-NumCell $wsShotsSummary 20 9 31.1
Set-TextCell $wsShotsSummary 21 4 "2025-11-06T19:30:00Z"
Set-TextCell $wsShotsSummary 22 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 22 5 37
Set-NumCell $wsShotsSummary 22 6 1012
Set-NumCell $wsShotsSummary 22 7 1295
Set-NumCell $wsShotsSummary 22 8 27.4
Set-NumCell $wsShotsSummary 22 9 35
Set-TextCell $wsShotsSummary 23 4 "2025-11-06T19:30:00Z"
Set-NumCell $wsShotsSummary 23 5 39
Set-NumCell $wsShotsSummary 23 6 935
Set-NumCell $wsShotsSummary 23 7 1116
Set-NumCell $wsShotsSummary 23 8 24
Set-NumCell $wsShotsSummary 23 9 28.6

# --- Meta_ext: refresh as_of_utc + bump build_version ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")

Set-TextCell $wsMeta 2 2 "2025-11-06T19:30:00Z"
Set-NumCell $wsMeta 2 4 49
